$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" header in H1, matching the style of the other header cells (copy G1's
# formatting since Excel COM's Range.Style setter expects a named style, not a
# direct-format copy) then overwrite the pasted text with the real value.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add the corresponding value in H2
$ws.Range("H2").Value = 1
